$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5814984262604241
$ws.Range("D2").Value = 0.0338715492648447
$ws.Range("E2").Value = 0.1800552420678652
$ws.Range("F2").Value = 1.021705754210004
$ws.Range("G2").Value = 0.8797214613288418
$ws.Range("H2").Value = 0.9098932397180306
$ws.Range("I2").Value = 1.13065978319408
$ws.Range("K2").Value = 0.4975660889531355
$ws.Range("L2").Value = 0.1959596293382901

$ws.Range("B3").Value = 0.5582198886308731
$ws.Range("D3").Value = 0.03334697646742057
$ws.Range("E3").Value = 0.1809042964432894
$ws.Range("F3").Value = 1.011847325405753
$ws.Range("G3").Value = 0.8708116068648053
$ws.Range("H3").Value = 0.9107309576123441
$ws.Range("I3").Value = 1.138733861339183
$ws.Range("K3").Value = 0.4339329784369284
$ws.Range("L3").Value = 0.1837846991509622

$ws.Range("B4").Value = 0.5441744397543573
$ws.Range("D4").Value = 0.03302135496411651
$ws.Range("E4").Value = 0.1814612083669731
$ws.Range("F4").Value = 1.006456857268425
$ws.Range("G4").Value = 0.8659602311750092
$ws.Range("H4").Value = 0.9117116594159427
$ws.Range("I4").Value = 1.14420799687101
$ws.Range("K4").Value = 0.3947947305211414
$ws.Range("L4").Value = 0.1764098340239855

$ws.Range("B5").Value = 0.5385135076139989
$ws.Range("D5").Value = 0.03288778474425058
$ws.Range("E5").Value = 0.1816971230942053
$ws.Range("F5").Value = 1.004426563507266
$ws.Range("G5").Value = 0.8641386061764251
$ws.Range("H5").Value = 0.9122285053622079
$ws.Range("I5").Value = 1.146568608212419
$ws.Range("K5").Value = 0.3788289930031112
$ws.Range("L5").Value = 0.1734298414599209

$ws.Range("B6").Value = 0.5375773137577085
$ws.Range("D6").Value = 0.032865552888655
$ws.Range("E6").Value = 0.1817368388086873
$ws.Range("F6").Value = 1.004099475205692
$ws.Range("G6").Value = 0.8638454983548627
$ws.Range("H6").Value = 0.9123214041297416
$ws.Range("I6").Value = 1.14696842744835
$ws.Range("K6").Value = 0.3761769027111086
$ws.Range("L6").Value = 0.1729365464490513

$ws.Range("B7").Value = 0.5440978400000631
$ws.Range("D7").Value = 0.03301955712503002
$ws.Range("E7").Value = 0.1814643536554188
$ws.Range("F7").Value = 1.006428802759075
$ws.Range("G7").Value = 0.8659350356048918
$ws.Range("H7").Value = 0.9117181553219353
$ws.Range("I7").Value = 1.144239307187693
$ws.Range("K7").Value = 0.3945794777051219
$ws.Range("L7").Value = 0.1763695422210674

$ws.Range("B8").Value = 0.5734208284291071
$ws.Range("D8").Value = 0.03369141610120252
$ws.Range("E8").Value = 0.1803406211810481
$ws.Range("F8").Value = 1.018168871939281
$ws.Range("G8").Value = 0.8765205713078785
$ws.Range("H8").Value = 0.9100852559915751
$ws.Range("I8").Value = 1.133336470033345
$ws.Range("K8").Value = 0.4756395294918434
$ws.Range("L8").Value = 0.1917408495551456

$ws.Range("B9").Value = 0.63287314271696
$ws.Range("D9").Value = 0.03498047869740617
$ws.Range("E9").Value = 0.1784184975145928
$ws.Range("F9").Value = 1.046463390940161
$ws.Range("G9").Value = 0.9022122193944142
$ws.Range("H9").Value = 0.9105872046647363
$ws.Range("I9").Value = 1.116057720001677
$ws.Range("K9").Value = 0.63405786666317
$ws.Range("L9").Value = 0.2226825384252322

$ws.Range("B10").Value = 0.6777264579944244
$ws.Range("D10").Value = 0.03590972967327133
$ws.Range("E10").Value = 0.1771767329962413
$ws.Range("F10").Value = 1.070489346735329
$ws.Range("G10").Value = 0.9241257484500522
$ws.Range("H10").Value = 0.9132209144325145
$ws.Range("I10").Value = 1.105867364871159
$ws.Range("K10").Value = 0.7501215893364588
$ws.Range("L10").Value = 0.2459055092674731

$ws.Range("B11").Value = 0.6983834146896584
$ws.Range("D11").Value = 0.03632850785223951
$ws.Range("E11").Value = 0.1766485744436586
$ws.Range("F11").Value = 1.082127835313088
$ws.Range("G11").Value = 0.934761145225508
$ws.Range("H11").Value = 0.9149124467360537
$ws.Range("I11").Value = 1.101776167339757
$ws.Range("K11").Value = 0.802852608475348
$ws.Range("L11").Value = 0.2565775306268563

$ws.Range("B12").Value = 0.7062416955180026
$ws.Range("D12").Value = 0.03648651203739917
$ws.Range("E12").Value = 0.1764538363746397
$ws.Range("F12").Value = 1.086637330736522
$ws.Range("G12").Value = 0.9388848549365605
$ws.Range("H12").Value = 0.9156240462240248
$ws.Range("I12").Value = 1.100305299196265
$ws.Range("K12").Value = 0.8228107056882834
$ws.Range("L12").Value = 0.2606342585261672

$ws.Range("B13").Value = 0.704547682619534
$ws.Range("D13").Value = 0.03645250889006135
$ws.Range("E13").Value = 0.1764955428500512
$ws.Range("F13").Value = 1.085661577661611
$ws.Range("G13").Value = 0.9379924497959848
$ws.Range("H13").Value = 0.9154676289061285
$ws.Range("I13").Value = 1.100618589600359
$ws.Range("K13").Value = 0.8185128227041503
$ws.Range("L13").Value = 0.2597598822987663

$ws.Range("B14").Value = 0.6990292023836844
$ws.Range("D14").Value = 0.03634151859864332
$ws.Range("E14").Value = 0.1766324478222532
$ws.Range("F14").Value = 1.082496783181995
$ws.Range("G14").Value = 0.9350984725472813
$ws.Range("H14").Value = 0.914969565895035
$ws.Range("I14").Value = 1.101653586780238
$ws.Range("K14").Value = 0.8044947755540477
$ws.Range("L14").Value = 0.2569109703178469

$ws.Range("B15").Value = 0.6956536414205914
$ws.Range("D15").Value = 0.0362734582410873
$ws.Range("E15").Value = 0.1767169911089899
$ws.Range("F15").Value = 1.080571580495317
$ws.Range("G15").Value = 0.9333383837544886
$ws.Range("H15").Value = 0.9146737439701269
$ws.Range("I15").Value = 1.102297762408071
$ws.Range("K15").Value = 0.7959070043963834
$ws.Range("L15").Value = 0.2551679425000088

$ws.Range("B16").Value = 0.6763815301775651
$ws.Range("D16").Value = 0.03588228124518622
$ws.Range("E16").Value = 0.1772119869634352
$ws.Range("F16").Value = 1.069743036295591
$ws.Range("G16").Value = 0.9234441549904915
$ws.Range("H16").Value = 0.9131203063937221
$ws.Range("I16").Value = 1.106145697498555
$ws.Range("K16").Value = 0.7466741144386617
$ws.Range("L16").Value = 0.2452102334652579

$ws.Range("B17").Value = 0.6646231921906463
$ws.Range("D17").Value = 0.03564128904742603
$ws.Range("E17").Value = 0.1775250451759605
$ws.Range("F17").Value = 1.063281873983613
$ws.Range("G17").Value = 0.9175454597391166
$ws.Range("H17").Value = 0.9122937740213928
$ws.Range("I17").Value = 1.108645790023026
$ws.Range("K17").Value = 0.7164539200341835
$ws.Range("L17").Value = 0.2391290818129193

$ws.Range("B18").Value = 0.6578839426995557
$ws.Range("D18").Value = 0.03550230624378514
$ws.Range("E18").Value = 0.1777085658527184
$ws.Range("F18").Value = 1.05963229686418
$ws.Range("G18").Value = 0.9142154252107701
$ws.Range("H18").Value = 0.9118648139985197
$ws.Range("I18").Value = 1.11013501764544
$ws.Range("K18").Value = 0.6990657834853664
$ws.Range("L18").Value = 0.2356415105079606

$ws.Range("B19").Value = 0.655606255705635
$ws.Range("D19").Value = 0.03545518576710904
$ws.Range("E19").Value = 0.177771297216152
$ws.Range("F19").Value = 1.058408061417722
$ws.Range("G19").Value = 0.9130986951830096
$ws.Range("H19").Value = 0.911727548622622
$ws.Range("I19").Value = 1.110648042434796
$ws.Range("K19").Value = 0.6931773919130535
$ws.Range("L19").Value = 0.2344624224322729

$ws.Range("B20").Value = 0.6658724236194189
$ws.Range("D20").Value = 0.03566698150033432
$ws.Range("E20").Value = 0.1774913618428533
$ws.Range("F20").Value = 1.06396276903142
$ws.Range("G20").Value = 0.9181668895698465
$ws.Range("H20").Value = 0.9123769527863885
$ws.Range("I20").Value = 1.108374347056312
$ws.Range("K20").Value = 0.7196715639622084
$ws.Range("L20").Value = 0.2397753809263889

$ws.Range("B21").Value = 0.70064914186446
$ws.Range("D21").Value = 0.03637413491159336
$ws.Range("E21").Value = 0.1765920927822795
$ws.Range("F21").Value = 1.08342358250988
$ws.Range("G21").Value = 0.9359458864294652
$ws.Range("H21").Value = 0.9151139299498539
$ws.Range("I21").Value = 1.101347455112602
$ws.Range("K21").Value = 0.8086124899813569
$ws.Range("L21").Value = 0.2577473453696939

$ws.Range("B22").Value = 0.7235870211726478
$ws.Range("D22").Value = 0.03683292642887182
$ws.Range("E22").Value = 0.176035043918247
$ws.Range("F22").Value = 1.096738495366992
$ws.Range("G22").Value = 0.9481270886602431
$ws.Range("H22").Value = 0.9173168985463462
$ws.Range("I22").Value = 1.097211872354109
$ws.Range("K22").Value = 0.8666820992188491
$ws.Range("L22").Value = 0.2695831900267081

$ws.Range("B23").Value = 0.7113256370229237
$ws.Range("D23").Value = 0.03658837308406504
$ws.Range("E23").Value = 0.1763295501812743
$ws.Range("F23").Value = 1.089577427813282
$ws.Range("G23").Value = 0.9415742242409522
$ws.Range("H23").Value = 0.9161032018849653
$ws.Range("I23").Value = 1.099377275218231
$ws.Range("K23").Value = 0.8356947222155213
$ws.Range("L23").Value = 0.2632579411477849

$ws.Range("B24").Value = 0.6653075813446776
$ws.Range("D24").Value = 0.03565536729023222
$ws.Range("E24").Value = 0.1775065790529911
$ws.Range("F24").Value = 1.063654733831115
$ws.Range("G24").Value = 0.9178857505262812
$ws.Range("H24").Value = 0.9123392036807019
$ws.Range("I24").Value = 1.108496904785731
$ws.Range("K24").Value = 0.7182169111796384
$ws.Range("L24").Value = 0.2394831624218341

$ws.Range("B25").Value = 0.6165825235524665
$ws.Range("D25").Value = 0.03463484942630402
$ws.Range("E25").Value = 0.1789084694384915
$ws.Range("F25").Value = 1.038242073425423
$ws.Range("G25").Value = 0.894730653139348
$ws.Range("H25").Value = 0.9100541649368381
$ws.Range("I25").Value = 1.120292568823203
$ws.Range("K25").Value = 0.5912590414043564
$ws.Range("L25").Value = 0.2142261603289199
